$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 136  # H33
$ws.Cells.Item(33, 9).Value = 216.4  # I33
$ws.Cells.Item(33, 11).Value = 216.4  # K33
$ws.Cells.Item(33, 13).Value = 12.59999999999999  # M33
$ws.Cells.Item(40, 8).Value = 2400  # H40
$ws.Cells.Item(40, 10).Value = 2400  # J40
$ws.Cells.Item(40, 12).Value = 2400  # L40
$ws.Cells.Item(40, 14).Value = -2750  # N40
$ws.Cells.Item(62, 8).Value = 11109.75  # H62
$ws.Cells.Item(62, 9).Value = 9998.333000000001  # I62
$ws.Cells.Item(62, 11).Value = 9998.333000000001  # K62
$ws.Cells.Item(62, 13).Value = -9374.333000000001  # M62
$ws.Cells.Item(65, 8).Value = 11109.75  # H65
$ws.Cells.Item(65, 9).Value = 9998.333000000001  # I65
$ws.Cells.Item(65, 11).Value = 49991.665  # K65
$ws.Cells.Item(65, 13).Value = -46871.665  # M65
$ws.Cells.Item(70, 8).Value = 2212  # H70
$ws.Cells.Item(70, 9).Value = 1750  # I70
$ws.Cells.Item(70, 10).Value = 2366  # J70
$ws.Cells.Item(70, 11).Value = 5250  # K70
$ws.Cells.Item(70, 12).Value = 7098  # L70
$ws.Cells.Item(70, 13).Value = -4980  # M70
$ws.Cells.Item(70, 14).Value = -7638  # N70
$ws.Cells.Item(73, 8).Value = 2212  # H73
$ws.Cells.Item(73, 9).Value = 1750  # I73
$ws.Cells.Item(73, 10).Value = 2366  # J73
$ws.Cells.Item(73, 11).Value = 5250  # K73
$ws.Cells.Item(73, 12).Value = 7098  # L73
$ws.Cells.Item(73, 13).Value = -4314  # M73
$ws.Cells.Item(73, 14).Value = -8970  # N73
$ws.Cells.Item(98, 8).Value = 994.1667  # H98
$ws.Cells.Item(98, 9).Value = 994.1667  # I98
$ws.Cells.Item(98, 10).Value = 0  # J98
$ws.Cells.Item(98, 11).Value = 994.1667  # K98
$ws.Cells.Item(98, 12).Value = 0  # L98
$ws.Cells.Item(98, 13).Value = 503.8333  # M98
$ws.Cells.Item(98, 14).ClearContents()  # N98
$ws.Cells.Item(111, 8).Value = 2130.7778  # H111
$ws.Cells.Item(111, 9).Value = 2166.3333  # I111
$ws.Cells.Item(111, 10).Value = 2059.6667  # J111
$ws.Cells.Item(111, 11).Value = 6498.999899999999  # K111
$ws.Cells.Item(111, 12).Value = 6179.000100000001  # L111
$ws.Cells.Item(111, 13).Value = -3431.999899999999  # M111
$ws.Cells.Item(111, 14).Value = -12313.0001  # N111
$ws.Cells.Item(122, 8).Value = 994.1667  # H122
$ws.Cells.Item(122, 9).Value = 994.1667  # I122
$ws.Cells.Item(122, 10).Value = 0  # J122
$ws.Cells.Item(122, 11).Value = 2982.5001  # K122
$ws.Cells.Item(122, 12).Value = 0  # L122
$ws.Cells.Item(122, 13).Value = -532.5001000000002  # M122
$ws.Cells.Item(122, 14).ClearContents()  # N122
$ws.Cells.Item(138, 8).Value = 3881.0557  # H138
$ws.Cells.Item(138, 9).Value = 4443.3335  # I138
$ws.Cells.Item(138, 10).Value = 3599.9167  # J138
$ws.Cells.Item(138, 11).Value = 13330.0005  # K138
$ws.Cells.Item(138, 12).Value = 10799.7501  # L138
$ws.Cells.Item(138, 13).Value = -8190.000499999998  # M138
$ws.Cells.Item(138, 14).Value = -21079.7501  # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2998.6  # H61
$ws.Cells.Item(61, 9).Value = 2333  # I61
$ws.Cells.Item(61, 11).Value = 2333  # K61
$ws.Cells.Item(61, 13).Value = -2121  # M61
$ws.Cells.Item(101, 8).Value = 0  # H101
$ws.Cells.Item(101, 10).Value = 0  # J101
$ws.Cells.Item(101, 12).Value = 0  # L101
$ws.Cells.Item(101, 14).ClearContents()  # N101
$ws.Cells.Item(112, 8).Value = 0  # H112
$ws.Cells.Item(112, 10).Value = 0  # J112
$ws.Cells.Item(112, 12).Value = 0  # L112
$ws.Cells.Item(112, 14).ClearContents()  # N112
$ws.Cells.Item(136, 8).Value = 2998.6  # H136
$ws.Cells.Item(136, 9).Value = 2333  # I136
$ws.Cells.Item(136, 11).Value = 6999  # K136
$ws.Cells.Item(136, 13).Value = -4449  # M136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(58, 8).Value = 41890  # H58
$ws.Cells.Item(58, 10).Value = 41890  # J58
$ws.Cells.Item(58, 12).Value = 41890  # L58
$ws.Cells.Item(58, 14).Value = -42478  # N58

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 34853984  # H86
$ws.Cells.Item(86, 9).Value = 34853984  # I86
$ws.Cells.Item(86, 11).Value = 34853984  # K86
$ws.Cells.Item(86, 13).Value = -34852861  # M86
$ws.Cells.Item(89, 8).Value = 34853984  # H89
$ws.Cells.Item(89, 9).Value = 34853984  # I89
$ws.Cells.Item(89, 11).Value = 174269920  # K89
$ws.Cells.Item(89, 13).Value = -174264304  # M89
$ws.Cells.Item(114, 8).Value = 93967.71000000001  # H114
$ws.Cells.Item(114, 10).Value = 93967.71000000001  # J114
$ws.Cells.Item(114, 12).Value = 93967.71000000001  # L114
$ws.Cells.Item(114, 14).Value = -102645.71  # N114
$ws.Cells.Item(131, 8).Value = 80000  # H131
$ws.Cells.Item(131, 10).Value = 80000  # J131
$ws.Cells.Item(131, 12).Value = 80000  # L131
$ws.Cells.Item(131, 14).Value = -90080  # N131
$ws.Cells.Item(141, 8).Value = 324376.25  # H141
$ws.Cells.Item(141, 10).Value = 324376.25  # J141
$ws.Cells.Item(141, 12).Value = 324376.25  # L141
$ws.Cells.Item(141, 14).Value = -334736.25  # N141

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 455.13333  # H23
$ws.Cells.Item(23, 9).Value = 268.8  # I23
$ws.Cells.Item(23, 10).Value = 548.3  # J23
$ws.Cells.Item(23, 11).Value = 806.4000000000001  # K23
$ws.Cells.Item(23, 12).Value = 1644.9  # L23
$ws.Cells.Item(23, 13).Value = -571.4000000000001  # M23
$ws.Cells.Item(23, 14).Value = -2114.9  # N23
$ws.Cells.Item(87, 8).Value = 0  # H87
$ws.Cells.Item(87, 9).Value = 0  # I87
$ws.Cells.Item(87, 11).Value = 0  # K87
$ws.Cells.Item(87, 13).ClearContents()  # M87
$ws.Cells.Item(90, 8).Value = 0  # H90
$ws.Cells.Item(90, 9).Value = 0  # I90
$ws.Cells.Item(90, 11).Value = 0  # K90
$ws.Cells.Item(90, 13).ClearContents()  # M90
$ws.Cells.Item(132, 8).Value = 724.5  # H132
$ws.Cells.Item(132, 9).Value = 724.5  # I132
$ws.Cells.Item(132, 11).Value = 6520.5  # K132
$ws.Cells.Item(132, 13).Value = -3990.5  # M132

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 642.5806  # H2
$ws.Cells.Item(2, 9).Value = 773.45  # I2
$ws.Cells.Item(2, 10).Value = 404.63635  # J2
$ws.Cells.Item(2, 11).Value = 773.45  # K2
$ws.Cells.Item(2, 12).Value = 404.63635  # L2
$ws.Cells.Item(2, 13).Value = -660.45  # M2
$ws.Cells.Item(2, 14).Value = -630.63635  # N2
$ws.Cells.Item(70, 8).Value = 4504  # H70
$ws.Cells.Item(70, 9).Value = 4008  # I70
$ws.Cells.Item(70, 10).Value = 5000  # J70
$ws.Cells.Item(70, 11).Value = 4008  # K70
$ws.Cells.Item(70, 12).Value = 5000  # L70
$ws.Cells.Item(70, 13).Value = -3738  # M70
$ws.Cells.Item(70, 14).Value = -5540  # N70
$ws.Cells.Item(73, 8).Value = 4504  # H73
$ws.Cells.Item(73, 9).Value = 4008  # I73
$ws.Cells.Item(73, 10).Value = 5000  # J73
$ws.Cells.Item(73, 11).Value = 4008  # K73
$ws.Cells.Item(73, 12).Value = 5000  # L73
$ws.Cells.Item(73, 13).Value = -3072  # M73
$ws.Cells.Item(73, 14).Value = -6872  # N73
$ws.Cells.Item(102, 8).Value = 1667.1875  # H102
$ws.Cells.Item(102, 9).Value = 1128.8462  # I102
$ws.Cells.Item(102, 11).Value = 1128.8462  # K102
$ws.Cells.Item(102, 13).Value = 493.1538  # M102
$ws.Cells.Item(107, 8).Value = 733.3077  # H107
$ws.Cells.Item(107, 9).Value = 536.8570999999999  # I107
$ws.Cells.Item(107, 10).Value = 962.5  # J107
$ws.Cells.Item(107, 11).Value = 536.8570999999999  # K107
$ws.Cells.Item(107, 12).Value = 962.5  # L107
$ws.Cells.Item(107, 13).Value = 1383.1429  # M107
$ws.Cells.Item(107, 14).Value = -4802.5  # N107
$ws.Cells.Item(111, 8).Value = 28700  # H111
$ws.Cells.Item(111, 10).Value = 28700  # J111
$ws.Cells.Item(111, 12).Value = 28700  # L111
$ws.Cells.Item(111, 14).Value = -34834  # N111
$ws.Cells.Item(122, 8).Value = 2080.5386  # H122
$ws.Cells.Item(122, 9).Value = 1299.875  # I122
$ws.Cells.Item(122, 10).Value = 3329.6  # J122
$ws.Cells.Item(122, 11).Value = 3899.625  # K122
$ws.Cells.Item(122, 12).Value = 9988.799999999999  # L122
$ws.Cells.Item(122, 13).Value = -1449.625  # M122
$ws.Cells.Item(122, 14).Value = -14888.8  # N122
$ws.Cells.Item(132, 8).Value = 3075.6924  # H132
$ws.Cells.Item(132, 9).Value = 2999.16  # I132
$ws.Cells.Item(132, 11).Value = 8997.48  # K132
$ws.Cells.Item(132, 13).Value = -6467.48  # M132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6448  # H40
$ws.Cells.Item(40, 9).Value = 6116.25  # I40
$ws.Cells.Item(40, 10).Value = 7332.6665  # J40
$ws.Cells.Item(40, 11).Value = 6116.25  # K40
$ws.Cells.Item(40, 12).Value = 7332.6665  # L40
$ws.Cells.Item(40, 13).Value = -5980.25  # M40
$ws.Cells.Item(40, 14).Value = -7604.6665  # N40
$ws.Cells.Item(110, 8).Value = 0  # H110
$ws.Cells.Item(110, 10).Value = 0  # J110
$ws.Cells.Item(110, 12).Value = 0  # L110
$ws.Cells.Item(110, 14).ClearContents()  # N110
$ws.Cells.Item(132, 8).Value = 6537  # H132
$ws.Cells.Item(132, 9).Value = 4895.3687  # I132
$ws.Cells.Item(132, 11).Value = 14686.1061  # K132
$ws.Cells.Item(132, 13).Value = -12156.1061  # M132
$ws.Cells.Item(136, 8).Value = 4015.84  # H136
$ws.Cells.Item(136, 9).Value = 3033.0833  # I136
$ws.Cells.Item(136, 11).Value = 9099.249899999999  # K136
$ws.Cells.Item(136, 13).Value = -6549.249899999999  # M136
